# Fruta / hortaliza, semanal
# Insert two new weekly records (rows 288-289) above the existing data,
# pushing all subsequent rows down by two (old row 288 -> new row 290, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at position 288 (existing rows 288..301 shift to 290..303)
$ws.Rows("288:289").Insert()

# --- New row 288: Betarraga, Primera ---
$ws.Cells.Item(288, 1).Value = 8
$ws.Cells.Item(288, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(288, 3).Value = "Coquimbo"
$ws.Cells.Item(288, 4).Value = 44747
$ws.Cells.Item(288, 5).Value = 4
$ws.Cells.Item(288, 6).Value = 100114014
$ws.Cells.Item(288, 7).Value = "Betarraga"
$ws.Cells.Item(288, 8).Value = "Sin especificar"
$ws.Cells.Item(288, 9).Value = "Primera"
$ws.Cells.Item(288, 10).Value = 2440
$ws.Cells.Item(288, 11).Value = 500
$ws.Cells.Item(288, 12).Value = 600
$ws.Cells.Item(288, 13).Value = 550
$ws.Cells.Item(288, 14).Value = "$/paquete 3 unidades"
$ws.Cells.Item(288, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(288, 16).Value = 183
$ws.Cells.Item(288, 17).Value = 3
$ws.Cells.Item(288, 18).Value = "Hortaliza"

# --- New row 289: Betarraga, Segunda ---
$ws.Cells.Item(289, 1).Value = 8
$ws.Cells.Item(289, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(289, 3).Value = "Coquimbo"
$ws.Cells.Item(289, 4).Value = 44747
$ws.Cells.Item(289, 5).Value = 4
$ws.Cells.Item(289, 6).Value = 100114014
$ws.Cells.Item(289, 7).Value = "Betarraga"
$ws.Cells.Item(289, 8).Value = "Sin especificar"
$ws.Cells.Item(289, 9).Value = "Segunda"
$ws.Cells.Item(289, 10).Value = 1560
$ws.Cells.Item(289, 11).Value = 400
$ws.Cells.Item(289, 12).Value = 450
$ws.Cells.Item(289, 13).Value = 425
$ws.Cells.Item(289, 14).Value = "$/paquete 3 unidades"
$ws.Cells.Item(289, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(289, 16).Value = 142
$ws.Cells.Item(289, 17).Value = 3
$ws.Cells.Item(289, 18).Value = "Hortaliza"
